$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 184, pushing existing rows 184-189 down to 186-191
$ws.Rows("184:185").Insert()

# New row 184 (Primera, week of 44448)
$ws.Cells.Item(184, 1).Value = 11
$ws.Cells.Item(184, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(184, 3).Value = "Bíobío"
$ws.Cells.Item(184, 4).Value2 = 44448
$ws.Cells.Item(184, 5).Value = 8
$ws.Cells.Item(184, 6).Value = 100112006
$ws.Cells.Item(184, 7).Value = "Repollo"
$ws.Cells.Item(184, 8).Value = "Crespo record"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 1000
$ws.Cells.Item(184, 11).Value = 700
$ws.Cells.Item(184, 12).Value = 800
$ws.Cells.Item(184, 13).Value = 750
$ws.Cells.Item(184, 14).Value = "$/unidad"
$ws.Cells.Item(184, 15).Value = "Región Metropolitana"
$ws.Cells.Item(184, 16).Value = 750
$ws.Cells.Item(184, 17).Value = 1
$ws.Cells.Item(184, 18).Value = "Hortaliza"

# New row 185 (Segunda, week of 44448)
$ws.Cells.Item(185, 1).Value = 11
$ws.Cells.Item(185, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(185, 3).Value = "Bíobío"
$ws.Cells.Item(185, 4).Value2 = 44448
$ws.Cells.Item(185, 5).Value = 8
$ws.Cells.Item(185, 6).Value = 100112006
$ws.Cells.Item(185, 7).Value = "Repollo"
$ws.Cells.Item(185, 8).Value = "Crespo record"
$ws.Cells.Item(185, 9).Value = "Segunda"
$ws.Cells.Item(185, 10).Value = 500
$ws.Cells.Item(185, 11).Value = 600
$ws.Cells.Item(185, 12).Value = 600
$ws.Cells.Item(185, 13).Value = 600
$ws.Cells.Item(185, 14).Value = "$/unidad"
$ws.Cells.Item(185, 15).Value = "Región Metropolitana"
$ws.Cells.Item(185, 16).Value = 600
$ws.Cells.Item(185, 17).Value = 1
$ws.Cells.Item(185, 18).Value = "Hortaliza"
